# Auto-generated Excel COM-interop script to apply scheduled-runner updates
# to currentAveragePrice* / LevePrice* / LeveProfit* columns across all 8 sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2066.3333
$ws.Range("I28").Value = 2066.3333
$ws.Range("K28").Value = 2066.3333
$ws.Range("M28").Value = -1581.3333
$ws.Range("H34").Value = 2948
$ws.Range("I34").Value = 2948
$ws.Range("K34").Value = 2948
$ws.Range("M34").Value = -2745
$ws.Range("H36").Value = 2948
$ws.Range("I36").Value = 2948
$ws.Range("K36").Value = 2948
$ws.Range("M36").Value = -2233
$ws.Range("H41").Value = 335.54544
$ws.Range("I41").Value = 271.22223
$ws.Range("J41").Value = 625
$ws.Range("K41").Value = 271.22223
$ws.Range("L41").Value = 625
$ws.Range("M41").Value = 168.77777
$ws.Range("N41").Value = -1505
$ws.Range("H53").Value = 169.54546
$ws.Range("I53").Value = 198.16667
$ws.Range("K53").Value = 198.16667
$ws.Range("M53").Value = 438.83333
$ws.Range("H86").Value = 3617
$ws.Range("I86").Value = 3617
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3617
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -2494
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 3617
$ws.Range("I89").Value = 3617
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 18085
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -12469
$ws.Range("N89").ClearContents()
$ws.Range("H93").Value = 57499.5
$ws.Range("J93").Value = 57499.5
$ws.Range("L93").Value = 57499.5
$ws.Range("N93").Value = -62491.5
$ws.Range("H105").Value = 10610
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H127").Value = 4250
$ws.Range("I127").Value = 2000
$ws.Range("K127").Value = 6000
$ws.Range("M127").Value = -1040
$ws.Range("H135").Value = 5046.8335
$ws.Range("I135").Value = 5046.8335
$ws.Range("K135").Value = 45421.5015
$ws.Range("M135").Value = -42886.5015
$ws.Range("H137").Value = 2563.8572
$ws.Range("I137").Value = 2490.4546
$ws.Range("J137").Value = 2833
$ws.Range("K137").Value = 7471.3638
$ws.Range("L137").Value = 8499
$ws.Range("M137").Value = -4921.3638
$ws.Range("N137").Value = -13599

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 44.42857
$ws.Range("I5").Value = 30.5
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 30.5
$ws.Range("L5").Value = 50
$ws.Range("M5").Value = 81.5
$ws.Range("N5").Value = -274
$ws.Range("H50").Value = 13558.444
$ws.Range("I50").Value = 4408.2
$ws.Range("K50").Value = 4408.2
$ws.Range("M50").Value = -3694.2
$ws.Range("H106").Value = 27332.666
$ws.Range("J106").Value = 27332.666
$ws.Range("L106").Value = 27332.666
$ws.Range("N106").Value = -29856.666
$ws.Range("H122").Value = 2649.7
$ws.Range("I122").Value = 2874.625
$ws.Range("K122").Value = 8623.875
$ws.Range("M122").Value = -6173.875

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 44.42857
$ws.Range("I4").Value = 30.5
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 30.5
$ws.Range("L4").Value = 50
$ws.Range("M4").Value = 84.5
$ws.Range("N4").Value = -280
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H22").Value = 1020.5714
$ws.Range("I22").Value = 1020.5714
$ws.Range("K22").Value = 1020.5714
$ws.Range("M22").Value = -847.5714
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H80").Value = 1531.3334
$ws.Range("I80").Value = 1531.3334
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 1531.3334
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -533.3334
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 1531.3334
$ws.Range("I83").Value = 1531.3334
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 7656.666999999999
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -2664.666999999999
$ws.Range("N83").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 5276.778
$ws.Range("I35").Value = 5279.2
$ws.Range("K35").Value = 5279.2
$ws.Range("M35").Value = -4985.2
$ws.Range("H107").Value = 759.5
$ws.Range("J107").Value = 1013
$ws.Range("L107").Value = 1013
$ws.Range("N107").Value = -4853

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 20000
$ws.Range("J57").Value = 20000
$ws.Range("L57").Value = 60000
$ws.Range("N57").Value = -61118
$ws.Range("H68").Value = 3119.95
$ws.Range("I68").Value = 2802
$ws.Range("J68").Value = 3136.6843
$ws.Range("K68").Value = 8406
$ws.Range("L68").Value = 9410.052899999999
$ws.Range("M68").Value = -7595
$ws.Range("N68").Value = -11032.0529
$ws.Range("H71").Value = 3119.95
$ws.Range("I71").Value = 2802
$ws.Range("J71").Value = 3136.6843
$ws.Range("K71").Value = 25218
$ws.Range("L71").Value = 28230.1587
$ws.Range("M71").Value = -21162
$ws.Range("N71").Value = -36342.1587

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 4833.3335
$ws.Range("J36").Value = 5250
$ws.Range("L36").Value = 5250
$ws.Range("N36").Value = -6220
$ws.Range("H101").Value = 24499.5
$ws.Range("J101").Value = 24499.5
$ws.Range("L101").Value = 24499.5
$ws.Range("N101").Value = -30989.5
$ws.Range("H122").Value = 4374.8335
$ws.Range("I122").Value = 4916.6665
$ws.Range("J122").Value = 3833
$ws.Range("K122").Value = 14749.9995
$ws.Range("L122").Value = 11499
$ws.Range("M122").Value = -12299.9995
$ws.Range("N122").Value = -16399
$ws.Range("H126").Value = 8538.9
$ws.Range("I126").Value = 8816.5
$ws.Range("J126").Value = 8122.5
$ws.Range("K126").Value = 26449.5
$ws.Range("L126").Value = 24367.5
$ws.Range("M126").Value = -23979.5
$ws.Range("N126").Value = -29307.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 855.44446
$ws.Range("I16").Value = 824.875
$ws.Range("K16").Value = 824.875
$ws.Range("M16").Value = -654.875
$ws.Range("H22").Value = 1736.1666
$ws.Range("J22").Value = 1736.1666
$ws.Range("L22").Value = 1736.1666
$ws.Range("N22").Value = -2326.1666
$ws.Range("H27").Value = 1736.1666
$ws.Range("J27").Value = 1736.1666
$ws.Range("L27").Value = 1736.1666
$ws.Range("N27").Value = -1950.1666
$ws.Range("H40").Value = 25714.857
$ws.Range("I40").Value = 7999.8
$ws.Range("K40").Value = 7999.8
$ws.Range("M40").Value = -7863.8
$ws.Range("H88").Value = 12000
$ws.Range("I88").Value = 12000
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 12000
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -11572
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 12000
$ws.Range("I91").Value = 12000
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 12000
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -10518
$ws.Range("N91").ClearContents()
$ws.Range("H97").Value = 22344
$ws.Range("J97").Value = 22344
$ws.Range("L97").Value = 22344
$ws.Range("N97").Value = -24326
$ws.Range("H103").Value = 46500
$ws.Range("J103").Value = 46500
$ws.Range("L103").Value = 46500
$ws.Range("N103").Value = -48844
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 26137.2
$ws.Range("I55").Value = 14690
$ws.Range("K55").Value = 14690
$ws.Range("M55").Value = -14413
$ws.Range("H104").Value = 17499.5
$ws.Range("J104").Value = 17499.5
$ws.Range("L104").Value = 17499.5
$ws.Range("N104").Value = -24487.5
$ws.Range("H126").Value = 1244.15
$ws.Range("I126").Value = 1232.1111
$ws.Range("K126").Value = 3696.3333
$ws.Range("M126").Value = -1226.3333
$ws.Range("H136").Value = 8860
$ws.Range("I136").Value = 7581.909
$ws.Range("J136").Value = 12374.75
$ws.Range("K136").Value = 22745.727
$ws.Range("L136").Value = 37124.25
$ws.Range("M136").Value = -20195.727
$ws.Range("N136").Value = -42224.25
